$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A, shifting the existing A:K header row (and the
# column-width definitions that go with it) one column to the right, to B:L.
$ws.Columns("A:A").Insert()

# Scroll the view back to the top-left (A1) and move the active selection.
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("E8").Select()
